$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Data Harian - Table")

# Add a new worksheet after the existing one and name it "Sheet1"
$newSheet = $wb.Worksheets.Add($null, $ws1)
$newSheet.Name = "Sheet1"

# Copy the finalized data table (header + 31 daily rows) and paste values only
$srcRange = $ws1.Range("A9:K40")
$srcRange.Copy()
$destRange = $newSheet.Range("A1:K32")
$destRange.PasteSpecial(-4163)
$excel.CutCopyMode = $false

# Make the new sheet the active/selected tab
$newSheet.Activate()
$newSheet.Range("A1:K32").Select()
